$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Clear "Clase secc 1" text from D7 and D8 (section 1 class removed), and reset
# their formatting to a plain blank cell (same style as other empty cells, e.g. C2/E9)
$ws.Range("D7").ClearContents()
$ws.Range("C2").Copy()
$ws.Range("D7").PasteSpecial(-4122)

$ws.Range("D8").ClearContents()
$ws.Range("C2").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update current selection to reflect the active cell at time of save
$ws.Range("G7").Select()
